$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data rows for August 29 (row 91, "raw" SSA data, styled/wrapped) and
# August 30 (row 92, "clean" data, plain formatting) appended below the
# existing historical log (which currently ends at row 90 / 2020-08-28).
# ---------------------------------------------------------------------------

# Write the date label for row 92 first, then row 91, so the shared-string
# table gets "2020-08-30" before "2020-08-29" (matches source ordering).
$ws.Range("A92").Formula = '="2020-08-30"'
$ws.Range("A92").Copy() | Out-Null
$ws.Range("A92").PasteSpecial(-4163) | Out-Null   # xlPasteValues -> keep as text, not a date serial

$ws.Range("A91").Formula = '="2020-08-29"'
$ws.Range("A91").Copy() | Out-Null
$ws.Range("A91").PasteSpecial(-4163) | Out-Null   # xlPasteValues -> keep as text, not a date serial

# Row 91 - raw data (Confirmados, Negativos, Sospechosos, Defunciones, %Hosp)
$ws.Range("B91").Value = 591712
$ws.Range("C91").Value = 658309
$ws.Range("D91").Value = 84310
$ws.Range("E91").Value = 63819
$ws.Range("F91").Value = 25.503285382077767

# Row 92 - clean data
$ws.Range("B92").Value = 595841
$ws.Range("C92").Value = 663474
$ws.Range("D92").Value = 81151
$ws.Range("E92").Value = 64158
$ws.Range("F92").Value = 25.46

# ---------------------------------------------------------------------------
# Formatting for row 91: 12pt font + wrapped text, taller row, and three
# extra (empty) formatted cells trailing the data (G91:I91).
# Build the combined format on a scratch cell first so the whole B91:I91
# range picks up a single new style instead of one per property change.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z1")
$scratch.Font.Size = 12
$scratch.WrapText = $true
$scratch.Copy() | Out-Null
$ws.Range("B91:I91").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$scratch.Clear() | Out-Null

$ws.Rows.Item(91).RowHeight = 16

# ---------------------------------------------------------------------------
# View state: scroll to show the new rows, select A91 (matches the source
# workbook's last-saved selection).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 68
$ws.Range("A91").Select() | Out-Null
